# Re-order the student roster: "Oscar" moves down below "Tom" and "Jay".
# Original rows (2-5): John, Oscar, Tom, Jay
# New rows      (2-5): John, Tom, Jay, Oscar
# Grades are unaffected (every student already has a 4), so only column A
# for rows 3-5 needs to be rewritten.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Tom"
$ws.Range("A4").Value = "Jay"
$ws.Range("A5").Value = "Oscar"
